{"js": "// Office.js (Word JavaScript API) script\n// Reproduces:\n//   1) document.xml: the \"Center-h4\" heading paragraph (styled Heading 4,\n//      centered) had its text split across two runs (\"Center-h\" + \"4\");\n//      it is retyped/normalized into a single run \"Center-h4\".\n//   2) styles.xml: Heading 4 gains bold; Heading 5 gains bold + italic\n//      and its size is bumped to 14pt (half-points 28).\n\n// --- 1) Merge the \"Center-h4\" Heading 4 paragraph into a single run ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nlet centerH4 = null;\nfor (const p of paragraphs.items) {\n  if (p.style === \"Heading 4\" && p.text === \"Center-h4\") {\n    centerH4 = p;\n    break;\n  }\n}\nif (!centerH4) {\n  throw new Error('Could not find the \"Center-h4\" Heading 4 paragraph.');\n}\n// Re-insert the same visible text as one continuous run, replacing the\n// two existing runs (\"Center-h\" and \"4\").\ncenterH4.getRange().insertText(\"Center-h4\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2) Update Heading 4 / Heading 5 style formatting ---\nconst heading4 = context.document.getStyles().getByNameOrNullObject(\"Heading 4\");\nconst heading5 = context.document.getStyles().getByNameOrNullObject(\"Heading 5\");\nawait context.sync();\n\nheading4.font.bold = true;\n\nheading5.font.bold = true;\nheading5.font.italic = true;\nheading5.font.size = 14;\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Reproduces:\n#   1) document.xml: the \"Center-h4\" heading paragraph (styled Heading 4,\n#      centered) had its text split across two runs (\"Center-h\" + \"4\");\n#      it is retyped/normalized into a single run \"Center-h4\".\n#   2) styles.xml: Heading 4 gains bold; Heading 5 gains bold + italic\n#      and its size is bumped to 14pt.\n\n$d = $word.ActiveDocument\n\n# --- 1) Merge the \"Center-h4\" Heading 4 paragraph into a single run ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Style = $d.Styles(\"Heading 4\")\n$find.Execute(\"Center-h4\", $false, $false, $false, $false, $false, $true, 1, $false, \"Center-h4\", 2)\n\n# --- 2) Update Heading 4 / Heading 5 style formatting ---\n$heading4 = $d.Styles(\"Heading 4\")\n$heading4.Font.Bold = $true\n\n$heading5 = $d.Styles(\"Heading 5\")\n$heading5.Font.Bold = $true\n$heading5.Font.Italic = $true\n$heading5.Font.Size = 14\n"}
